$wb = $excel.ActiveWorkbook

$wsQuery    = $wb.Worksheets.Item("Query")
$wsPrepared = $wb.Worksheets.Item("Prepared")
$wsLess     = $wb.Worksheets.Item("LessGreater")

# --- Ticket 41 part 2: "replacing expressions with values, continue replacing
# them only one at a time" - the SQL literal in the LessGreater sheet is
# reformatted onto multiple lines (embedded newlines inside the single-quoted
# SQL string) -----------------------------------------------------------
$newSql = "<jt:forEach items=""`${jdbc.execQuery('SELECT *`nFROM employee`nWHERE first_name <> \'Randy\'')}"" var=""employee"" >`${employee.first_name}"
$wsLess.Range("A2").Value = $newSql

# Give that same cell a (new, distinct) cell style - mirrors the extra
# cellXfs entry added to styles.xml for this cell.
$wsLess.Range("A2").HorizontalAlignment = 1  # xlGeneral

# --- Selection / active-tab bookkeeping -------------------------------
# Prepared: selection becomes the header row A1:G1 (no single active cell)
$wsPrepared.Range("A1:G1").Select()

# LessGreater should no longer be the active/selected tab.
# Query becomes the active sheet/tab instead.
$wsQuery.Activate()
